$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the two "UNT" placeholder values to distinct UNT1 / UNT2 values
$ws.Range("B4").Value = "UNT1"
$ws.Range("C4").Value = "UNT2"

# Move the active selection from B2 to B1
$ws.Range("B1").Select()
